$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = "JH55X7"
$ws.Range("B44").Value = "Estación o plancha para soldar y desoldar LED TV"
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 120000
$ws.Range("F44").Value = 3
$ws.Range("G44").Value = 4
$ws.Range("H44").Formula = "=(E44-D44)*G44"
$ws.Range("I44").Formula = "=D44*F44"
$ws.Range("J44").Value = 0
